$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the MSME table (row 9) - bold, matches existing "title" style
$ws.Range("B9").Value = "Number of employees"
$ws.Range("C9").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D9").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B9:D9").Font.Bold = $true

# Data rows 10-13 (Micro / Small / Medium / Large enterprise classes).
# The Assets/Turnover columns are blank in the source data, but still need
# to be real (empty) text cells rather than no cell at all, so a lone
# apostrophe (Excel's "force text" prefix) is used to get an empty string
# value instead of clearing the cell; the style is then reset back to
# Normal so the quote-prefix formatting doesn't stick around.
$ws.Range("A10").Value = "Micro"
$ws.Range("B10").Value = "0-3"
$ws.Range("C10").Value = "'"
$ws.Range("D10").Value = "'"
$ws.Range("C10:D10").Style = "Normal"

$ws.Range("A11").Value = "Small"
$ws.Range("B11").Value = "4-20"
$ws.Range("C11").Value = "'"
$ws.Range("D11").Value = "'"
$ws.Range("C11:D11").Style = "Normal"

$ws.Range("A12").Value = "Medium"
$ws.Range("B12").Value = "21-50"
$ws.Range("C12").Value = "'"
$ws.Range("D12").Value = "'"
$ws.Range("C12:D12").Style = "Normal"

$ws.Range("A13").Value = "Large"
$ws.Range("B13").Value = ">50"
$ws.Range("C13").Value = "'"
$ws.Range("D13").Value = "'"
$ws.Range("C13:D13").Style = "Normal"
